$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Cells.Item(2, 4)
$cell.Value = "'" + '70.348.05'
$cell.Style = "Normal"

$cell = $ws.Cells.Item(3, 4)
$cell.Value = "'" + '3.623.64'
$cell.Style = "Normal"
$ws.Cells.Item(3, 5).Value = '  +2.84%  '

$ws.Cells.Item(4, 5).Value = '  +0.01%  '

$cell = $ws.Cells.Item(5, 4)
$cell.Value = "'" + '601.53'
$cell.Style = "Normal"
$ws.Cells.Item(5, 5).Value = '  -0.87%  '

$cell = $ws.Cells.Item(6, 4)
$cell.Value = "'" + '196.80'
$cell.Style = "Normal"
$ws.Cells.Item(6, 5).Value = '  +0.03%  '

$cell = $ws.Cells.Item(7, 4)
$cell.Value = "'" + '0.625'
$cell.Style = "Normal"
$ws.Cells.Item(7, 5).Value = '  -0.87%  '

$cell = $ws.Cells.Item(8, 4)
$cell.Value = "'" + '1.00'
$cell.Style = "Normal"
$ws.Cells.Item(8, 5).Value = '  +0.09%  '

$cell = $ws.Cells.Item(9, 4)
$cell.Value = "'" + '0.212'
$cell.Style = "Normal"
$ws.Cells.Item(9, 5).Value = '  +6.50%  '

$ws.Cells.Item(10, 5).Value = '  -0.50%  '

$cell = $ws.Cells.Item(11, 4)
$cell.Value = "'" + '53.23'
$cell.Style = "Normal"
$ws.Cells.Item(11, 5).Value = '  -1.03%  '

$ws.Cells.Item(12, 5).Value = '  +0.73%  '

$cell = $ws.Cells.Item(13, 4)
$cell.Value = "'" + '9.54'
$cell.Style = "Normal"
$ws.Cells.Item(13, 5).Value = '  +0.44%  '

$cell = $ws.Cells.Item(14, 4)
$cell.Value = "'" + '4.197.31'
$cell.Style = "Normal"
$ws.Cells.Item(14, 5).Value = '  +2.80%  '

$cell = $ws.Cells.Item(15, 4)
$cell.Value = "'" + '604.77'
$cell.Style = "Normal"
$ws.Cells.Item(15, 5).Value = '  +0.98%  '

$cell = $ws.Cells.Item(16, 4)
$cell.Value = "'" + '12.96'
$cell.Style = "Normal"
$ws.Cells.Item(16, 5).Value = '  +1.35%  '

$cell = $ws.Cells.Item(17, 4)
$cell.Value = "'" + '70.459.42'
$cell.Style = "Normal"
$ws.Cells.Item(17, 5).Value = '  +0.63%  '

$ws.Cells.Item(18, 2).Value = 'WrappedEther'
$ws.Cells.Item(18, 3).Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$cell = $ws.Cells.Item(18, 4)
$cell.Value = "'" + '3.632.96'
$cell.Style = "Normal"
$ws.Cells.Item(18, 5).Value = '  +3.52%  '

$ws.Cells.Item(19, 2).Value = 'Chainlink'
$ws.Cells.Item(19, 3).Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$cell = $ws.Cells.Item(19, 4)
$cell.Value = "'" + '19.04'
$cell.Style = "Normal"
$ws.Cells.Item(19, 5).Value = '  -0.25%  '

$ws.Cells.Item(20, 5).Value = '  +1.34%  '

$ws.Cells.Item(21, 5).Value = '  +0.59%  '

$cell = $ws.Cells.Item(22, 4)
$cell.Value = "'" + '18.09'
$cell.Style = "Normal"
$ws.Cells.Item(22, 5).Value = '  -0.79%  '

$ws.Cells.Item(23, 5).Value = '  -0.97%  '

$cell = $ws.Cells.Item(24, 4)
$cell.Value = "'" + '103.35'
$cell.Style = "Normal"
$ws.Cells.Item(24, 5).Value = '  +1.10%  '

$ws.Cells.Item(25, 5).Value = '  -1.32%  '

$ws.Cells.Item(26, 5).Value = '  -6.47%  '

$cell = $ws.Cells.Item(27, 4)
$cell.Value = "'" + '10.59'
$cell.Style = "Normal"
$ws.Cells.Item(27, 5).Value = '  -2.67%  '

$cell = $ws.Cells.Item(28, 4)
$cell.Value = "'" + '9.68'
$cell.Style = "Normal"
$ws.Cells.Item(28, 5).Value = '  +0.70%  '

$cell = $ws.Cells.Item(29, 4)
$cell.Value = "'" + '33.83'
$cell.Style = "Normal"
$ws.Cells.Item(29, 5).Value = '  +1.32%  '

$cell = $ws.Cells.Item(30, 4)
$cell.Value = "'" + '4.67'
$cell.Style = "Normal"
$ws.Cells.Item(30, 5).Value = '  +7.87%  '

$cell = $ws.Cells.Item(31, 4)
$cell.Value = "'" + '7.24'
$cell.Style = "Normal"
$ws.Cells.Item(31, 5).Value = '  +2.14%  '

$ws.Cells.Item(32, 5).Value = '  -1.36%  '

$ws.Cells.Item(33, 5).Value = '  +0.72%  '

$cell = $ws.Cells.Item(34, 4)
$cell.Value = "'" + '63.23'
$cell.Style = "Normal"
$ws.Cells.Item(34, 5).Value = '  +0.17%  '

$cell = $ws.Cells.Item(35, 4)
$cell.Value = "'" + '0.0₃0886'
$cell.Style = "Normal"
$ws.Cells.Item(35, 5).Value = '  +3.86%  '

$cell = $ws.Cells.Item(36, 4)
$cell.Value = "'" + '3.942.60'
$cell.Style = "Normal"
$ws.Cells.Item(36, 5).Value = '  +5.49%  '

$ws.Cells.Item(37, 5).Value = '  +0.23%  '

$ws.Cells.Item(38, 2).Value = 'Fetch.AI'
$ws.Cells.Item(38, 3).Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$cell = $ws.Cells.Item(38, 4)
$cell.Value = "'" + '3.06'
$cell.Style = "Normal"
$ws.Cells.Item(38, 5).Value = '  -0.68%  '

$ws.Cells.Item(39, 2).Value = 'Bittensor'
$ws.Cells.Item(39, 3).Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$cell = $ws.Cells.Item(39, 4)
$cell.Value = "'" + '517.86'
$cell.Style = "Normal"
$ws.Cells.Item(39, 5).Value = '  +6.01%  '

$cell = $ws.Cells.Item(40, 4)
$cell.Value = "'" + '36.67'
$cell.Style = "Normal"
$ws.Cells.Item(40, 5).Value = '  +0.12%  '

$ws.Cells.Item(41, 5).Value = '  -1.15%  '

$ws.Cells.Item(42, 5).Value = '  -2.61%  '

$cell = $ws.Cells.Item(43, 4)
$cell.Value = "'" + '0.137'
$cell.Style = "Normal"
$ws.Cells.Item(43, 5).Value = '  +2.39%  '

$cell = $ws.Cells.Item(44, 4)
$cell.Value = "'" + '0.0460'
$cell.Style = "Normal"
$ws.Cells.Item(44, 5).Value = '  +1.32%  '

$ws.Cells.Item(45, 5).Value = '  +6.47%  '

$cell = $ws.Cells.Item(46, 4)
$cell.Value = "'" + '2.90'
$cell.Style = "Normal"
$ws.Cells.Item(46, 5).Value = '  +2.13%  '

$ws.Cells.Item(47, 5).Value = '  -0.06%  '

$ws.Cells.Item(48, 5).Value = '  +0.34%  '

$ws.Cells.Item(49, 5).Value = '  -0.23%  '

$cell = $ws.Cells.Item(50, 4)
$cell.Value = "'" + '0.000250'
$cell.Style = "Normal"
$ws.Cells.Item(50, 5).Value = '  +0.46%  '

$ws.Cells.Item(51, 5).Value = '  +0.62%  '
